$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 278.42856
$ws.Range("I19").Value = 327.2
$ws.Range("J19").Value = 156.5
$ws.Range("K19").Value = 327.2
$ws.Range("L19").Value = 156.5
$ws.Range("M19").Value = -152.2
$ws.Range("N19").Value = -506.5
$ws.Range("H40").Value = 2199.25
$ws.Range("J40").Value = 3699.5
$ws.Range("L40").Value = 3699.5
$ws.Range("N40").Value = -4049.5
$ws.Range("H51").Value = 5666.3335
$ws.Range("J51").Value = 5499.5
$ws.Range("L51").Value = 5499.5
$ws.Range("N51").Value = -6467.5
$ws.Range("H74").Value = 3543.3333
$ws.Range("I74").Value = 315
$ws.Range("K74").Value = 315
$ws.Range("M74").Value = 621
$ws.Range("H77").Value = 3543.3333
$ws.Range("I77").Value = 315
$ws.Range("K77").Value = 1575
$ws.Range("M77").Value = 3105
$ws.Range("H98").Value = 1663
$ws.Range("I98").Value = 1663
$ws.Range("K98").Value = 1663
$ws.Range("M98").Value = -165
$ws.Range("H107").Value = 60585.734
$ws.Range("I107").Value = 82471.37
$ws.Range("K107").Value = 82471.37
$ws.Range("M107").Value = -80551.37
$ws.Range("H122").Value = 1663
$ws.Range("I122").Value = 1663
$ws.Range("K122").Value = 4989
$ws.Range("M122").Value = -2539
$ws.Range("H129").Value = 997
$ws.Range("I129").Value = 997
$ws.Range("K129").Value = 2991
$ws.Range("M129").Value = 2009
$ws.Range("H131").Value = 3322.6667
$ws.Range("I131").Value = 499.5
$ws.Range("K131").Value = 1498.5
$ws.Range("M131").Value = 3541.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1900
$ws.Range("I63").Value = 1900
$ws.Range("K63").Value = 1900
$ws.Range("M63").Value = -1214
$ws.Range("H66").Value = 1900
$ws.Range("I66").Value = 1900
$ws.Range("K66").Value = 9500
$ws.Range("M66").Value = -6068

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2091.8333
$ws.Range("I5").Value = 4000
$ws.Range("J5").Value = 1710.2
$ws.Range("K5").Value = 4000
$ws.Range("L5").Value = 1710.2
$ws.Range("M5").Value = -3887
$ws.Range("N5").Value = -1936.2
$ws.Range("H11").Value = 420.66666
$ws.Range("I11").Value = 47.5
$ws.Range("J11").Value = 607.25
$ws.Range("K11").Value = 47.5
$ws.Range("L11").Value = 607.25
$ws.Range("M11").Value = 92.5
$ws.Range("N11").Value = -887.25
$ws.Range("H109").Value = 69000
$ws.Range("I109").Value = 69000
$ws.Range("K109").Value = 69000
$ws.Range("M109").Value = -67613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1130.9286
$ws.Range("J16").Value = 1992.5
$ws.Range("L16").Value = 1992.5
$ws.Range("N16").Value = -2566.5
$ws.Range("H113").Value = 1130.9286
$ws.Range("J113").Value = 1992.5
$ws.Range("L113").Value = 1992.5
$ws.Range("N113").Value = -6332.5
$ws.Range("H132").Value = 2197.25
$ws.Range("I132").Value = 2197.25
$ws.Range("K132").Value = 6591.75
$ws.Range("M132").Value = -4061.75
$ws.Range("H141").Value = 123471.71
$ws.Range("J141").Value = 123471.71
$ws.Range("L141").Value = 123471.71
$ws.Range("N141").Value = -133831.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 108.111115
$ws.Range("I23").Value = 65.166664
$ws.Range("J23").Value = 129.58333
$ws.Range("K23").Value = 195.499992
$ws.Range("L23").Value = 388.74999
$ws.Range("M23").Value = 39.50000800000001
$ws.Range("N23").Value = -858.74999
$ws.Range("H117").Value = 2263.8333
$ws.Range("J117").Value = 3684
$ws.Range("L117").Value = 11052
$ws.Range("N117").Value = -17936
$ws.Range("H121").Value = 662.6667
$ws.Range("J121").Value = 662.6667
$ws.Range("L121").Value = 1988.0001
$ws.Range("N121").Value = -4608.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 647.5
$ws.Range("I113").Value = 647.5
$ws.Range("K113").Value = 647.5
$ws.Range("M113").Value = 1522.5
$ws.Range("H122").Value = 5399.294
$ws.Range("I122").Value = 3816.2222
$ws.Range("J122").Value = 7180.25
$ws.Range("K122").Value = 11448.6666
$ws.Range("L122").Value = 21540.75
$ws.Range("M122").Value = -8998.6666
$ws.Range("N122").Value = -26440.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10851.137
$ws.Range("I7").Value = 10415.477
$ws.Range("J7").Value = 20000
$ws.Range("K7").Value = 10415.477
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = -10303.477
$ws.Range("N7").Value = -20224
$ws.Range("H22").Value = 1396.1428
$ws.Range("I22").Value = 982.75
$ws.Range("J22").Value = 1947.3334
$ws.Range("K22").Value = 982.75
$ws.Range("L22").Value = 1947.3334
$ws.Range("M22").Value = -687.75
$ws.Range("N22").Value = -2537.3334
$ws.Range("H27").Value = 1396.1428
$ws.Range("I27").Value = 982.75
$ws.Range("J27").Value = 1947.3334
$ws.Range("K27").Value = 982.75
$ws.Range("L27").Value = 1947.3334
$ws.Range("M27").Value = -875.75
$ws.Range("N27").Value = -2161.3334
$ws.Range("H38").Value = 33333
$ws.Range("J38").Value = 33333
$ws.Range("L38").Value = 33333
$ws.Range("N38").Value = -34153
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1000
$ws.Range("N46").ClearContents()
$ws.Range("M46").Value = -812
$ws.Range("H93").Value = 10419498
$ws.Range("I93").Value = 10419498
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 10419498
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H126").Value = 10851.137
$ws.Range("I126").Value = 10415.477
$ws.Range("J126").Value = 20000
$ws.Range("K126").Value = 31246.431
$ws.Range("L126").Value = 60000
$ws.Range("M126").Value = -28776.431
$ws.Range("N126").Value = -64940
$ws.Range("H132").Value = 5277.6665
$ws.Range("I132").Value = 4812.5
$ws.Range("K132").Value = 14437.5
$ws.Range("M132").Value = -11907.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 60000
$ws.Range("I15").Value = 60000
$ws.Range("K15").Value = 60000
$ws.Range("M15").Value = -59712
